$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 26 ("RM 232") entirely - remaining rows shift up.
$ws.Range("A26:F26").EntireRow.Delete()

# After the above delete, the former row 28 ("SC 92") is now row 27.
# Delete it too so everything shifts up again.
$ws.Range("A27:F27").EntireRow.Delete()

# The former row 34 ("SC 193") is now row 32; its column C value
# should be cleared (was 10.5, becomes blank).
$ws.Range("C32").ClearContents()
